$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.7252663449146157
$ws.Range("C2").Value = -0.7264947710836137
$ws.Range("D2").Value = -0.6438116996505815

$ws.Range("B3").Value = 0.7857583742924954
$ws.Range("C3").Value = 0.6945631837384654
$ws.Range("D3").Value = -0.5848874524283162

$ws.Range("B4").Value = 0.7734892888043006
$ws.Range("C4").Value = -0.6814466276472175
$ws.Range("D4").Value = 0.6447219063213621

$ws.Range("B5").Value = -0.6882343950236557
$ws.Range("C5").Value = 0.6013359868225709
$ws.Range("D5").Value = -0.492900360568737

$ws.Range("B6").Value = -0.8232155214576399
$ws.Range("C6").Value = -0.7364419089536182
$ws.Range("D6").Value = -0.5691473661445998

$ws.Range("B7").Value = -0.7997800534804209
$ws.Range("C7").Value = -0.4000529195126185
$ws.Range("D7").Value = 0.7873902526128235

$ws.Range("B8").Value = -0.8252920169842061
$ws.Range("C8").Value = -0.8089826795689997
$ws.Range("D8").Value = 0.688415227476193

$ws.Range("B9").Value = 0.7795322195061662
$ws.Range("C9").Value = 0.7280440006258964
$ws.Range("D9").Value = -0.6123770258533671
